# INDIANA_2017.xlsx cleanup:
#  1. Rename header columns to snake_case field names.
#  2. Title-case the lowercase Spanish connector words ("de", "del", "el",
#     "la", "las", "los", "y") that appear strictly between two other words
#     inside state / municipality names (columns A and B), e.g.
#     "Pabellón de Arteaga" -> "Pabellón De Arteaga". Leading/trailing
#     occurrences (e.g. "El Llano", "La Paz", "Del Nayar") are left alone
#     because they are not surrounded by a word on both sides.
#  3. Correct a floating point rounding difference in D265.
#  4. Drop the trailing metadata/footnote rows (1364-1368) so the used
#     range shrinks back down to A1:D1362.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B ---------------------
$words = @("de", "del", "el", "la", "las", "los", "y")
$lastRow = 1362

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            $parts = $v.Split(" ")
            $changed = $false
            for ($i = 1; $i -lt ($parts.Count - 1); $i++) {
                $p = $parts[$i]
                if ($words -contains $p) {
                    $parts[$i] = $p.Substring(0, 1).ToUpper() + $p.Substring(1)
                    $changed = $true
                }
            }
            if ($changed) {
                $cell.Value2 = [string]::Join(" ", $parts)
            }
        }
    }
}

# --- 3. Floating point fix -------------------------------------------------
$ws.Range("D265").Value = 0.09525299131930869

# --- 4. Drop trailing metadata rows ----------------------------------------
$ws.Range("A1363:D1368").ClearContents()
